# "Fruta / hortaliza, semanal" — add this week's new weekly price entry.
# A new data row is inserted right after the current row 40 (i.e. at row 41),
# pushing all the existing data rows (old 41..147) down by one (new 42..148).
# The new row 41 carries the new weekly record: same market/product
# attributes as the (now shifted) row 42, but with an updated date (D) and
# volume (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 41; everything below shifts down one row.
$ws.Rows("41").Insert()

# Populate the new row 41 with the new weekly record.
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 45177
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112010
$ws.Range("G41").Value = "Achicoria"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 100
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = 10000
$ws.Range("N41").Value = "`$/caja 18 unidades"
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 556
$ws.Range("Q41").Value = 18
$ws.Range("R41").Value = "Hortaliza"

# Make sure the new date cell keeps the date number format used by the
# rest of column D (style index 2 in the original workbook).
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
